$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.620.22"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").Value = "1.793.03"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'226.65"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("E6").Value = "  +1.74%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'32.81"
$ws.Range("E8").Value = "  +3.02%  "
$ws.Range("D9").Value = "'0.296"
$ws.Range("E9").Value = "  +1.77%  "
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").Value = "2.053.19"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").Value = "1.790.83"
$ws.Range("E14").Value = "  +0.25%  "
$ws.Range("D15").Value = "'0.635"
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").Value = "34.558.22"
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("D18").Value = "'68.74"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").Value = "'247.94"
$ws.Range("E19").Value = "  +0.76%  "
$ws.Range("E20").Value = "  +2.41%  "
$ws.Range("E21").Value = "  +2.36%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  +1.80%  "
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("D25").Value = "'165.53"
$ws.Range("E25").Value = "  +2.37%  "
$ws.Range("D26").Value = "'7.27"
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("D27").Value = "'16.55"
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("E28").Value = "  +2.13%  "
$ws.Range("D29").Value = "'1.01"
$ws.Range("E29").Value = "  +1.10%  "
$ws.Range("D30").Value = "'4.16"
$ws.Range("E30").Value = "  +14.02%  "
$ws.Range("D31").Value = "'3.82"
$ws.Range("E31").Value = "  +2.87%  "
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  +2.07%  "
$ws.Range("D35").Value = "1.427.23"
$ws.Range("E35").Value = "  -1.27%  "
$ws.Range("D36").Value = "'2.59"
$ws.Range("E36").Value = "  +6.43%  "
$ws.Range("D37").Value = "'0.671"
$ws.Range("E37").Value = "  +2.28%  "
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("E39").Value = "  +0.35%  "
$ws.Range("D40").Value = "'85.16"
$ws.Range("E40").Value = "  +6.21%  "
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("D42").Value = "'2.75"
$ws.Range("E42").Value = "  +2.73%  "
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").Value = "'13.58"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("E45").Value = "  +3.90%  "
$ws.Range("E46").Value = "  +1.01%  "
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("D48").Value = "1.953.40"
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("D49").Value = "'106.10"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("E50").Value = "  +0.05%  "

Write-Output "done"
